$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("2021-07-15", "overview", "K02000001", "United Kingdom", 5281098, 48553, 63, 128593),
    @("2021-07-16", "overview", "K02000001", "United Kingdom", 5332371, 51870, 49, 128642),
    @("2021-07-17", "overview", "K02000001", "United Kingdom", 5386340, 54674, 41, 128683),
    @("2021-07-18", "overview", "K02000001", "United Kingdom", 5433939, 48161, 25, 128708)
)

$startRow = 338
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 1).Style = "Normal"
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
}
